$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet logs one row per day (col A = date, B/C/D = counters).
# Last existing row is 366 (2021-09-01, serial 44440).
# "aggiornamento a 9/09 compreso" -> append rows through 2021-09-09
# (serials 44441..44448), matching the formatting of the last row.

$lastRow = 366
$newRowsCount = 8
$startSerial = 44441

# Copy the formatting (style) of the last data row down onto the new rows
# so the new date cells keep the same date number-format / style index
# instead of Excel inventing a brand-new style.
$srcRange = $ws.Range("A" + $lastRow + ":D" + $lastRow)
$srcRange.Copy() | Out-Null

$destRange = $ws.Range("A" + ($lastRow + 1) + ":D" + ($lastRow + $newRowsCount))
$destRange.PasteSpecial(-4122) | Out-Null

for ($i = 0; $i -lt $newRowsCount; $i++) {
    $row = $lastRow + 1 + $i
    $serial = $startSerial + $i

    $ws.Cells.Item($row, 1).Value = $serial
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
}
